# Scheduled market-data refresh: update computed profit columns (H..N)
# on the per-job sheets of Mateus_Profits.xlsx with freshly pulled prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 6570.56
$ws.Range("I15").Value = 6570.56
$ws.Range("K15").Value = 19711.68
$ws.Range("M15").Value = -19542.68

# Row 19
$ws.Range("H19").Value = 2428.25
$ws.Range("I19").Value = 1749.25
$ws.Range("K19").Value = 1749.25
$ws.Range("M19").Value = -1574.25

# Row 32
$ws.Range("H32").Value = 9594.857
$ws.Range("I32").Value = 5432.8
$ws.Range("K32").Value = 5432.8
$ws.Range("M32").Value = -5106.8

# Row 33
$ws.Range("H33").Value = 468.9524
$ws.Range("I33").Value = 454.92856
$ws.Range("J33").Value = 497
$ws.Range("K33").Value = 454.92856
$ws.Range("L33").Value = 497
$ws.Range("M33").Value = -225.92856
$ws.Range("N33").Value = -955

# Row 43
$ws.Range("H43").Value = 3299.6
$ws.Range("J43").Value = 3574.5
$ws.Range("L43").Value = 3574.5
$ws.Range("N43").Value = -3712.5

# Row 74
$ws.Range("H74").Value = 6334.1665
$ws.Range("I74").Value = 3689.0625
$ws.Range("J74").Value = 11624.375
$ws.Range("K74").Value = 3689.0625
$ws.Range("L74").Value = 11624.375
$ws.Range("M74").Value = -2753.0625
$ws.Range("N74").Value = -13496.375

# Row 77
$ws.Range("H77").Value = 6334.1665
$ws.Range("I77").Value = 3689.0625
$ws.Range("J77").Value = 11624.375
$ws.Range("K77").Value = 18445.3125
$ws.Range("L77").Value = 58121.875
$ws.Range("M77").Value = -13765.3125
$ws.Range("N77").Value = -67481.875

# Row 87
$ws.Range("H87").Value = 57500
$ws.Range("J87").Value = 57500
$ws.Range("L87").Value = 57500
$ws.Range("N87").Value = -59996

# Row 90
$ws.Range("H90").Value = 57500
$ws.Range("J90").Value = 57500
$ws.Range("L90").Value = 172500
$ws.Range("N90").Value = -184980

# Row 103
$ws.Range("H103").Value = 500
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()

# Row 132
$ws.Range("H132").Value = 8499.333000000001
$ws.Range("I132").Value = 9874.333000000001
$ws.Range("J132").Value = 2999.3333
$ws.Range("K132").Value = 29622.999
$ws.Range("L132").Value = 8997.999899999999
$ws.Range("M132").Value = -27092.999
$ws.Range("N132").Value = -14057.9999

# Row 141
$ws.Range("H141").Value = 5445.7617
$ws.Range("I141").Value = 4689.5557
$ws.Range("K141").Value = 14068.6671
$ws.Range("M141").Value = -8888.667099999999


$ws = $wb.Worksheets.Item("ARM")
# Row 55
$ws.Range("H55").Value = 34000
$ws.Range("J55").Value = 34500
$ws.Range("L55").Value = 34500
$ws.Range("N55").Value = -35130

# Row 88
$ws.Range("H88").Value = 1661.2307
$ws.Range("J88").Value = 1858.2858
$ws.Range("L88").Value = 1858.2858
$ws.Range("N88").Value = -2670.2858

# Row 91
$ws.Range("H91").Value = 1661.2307
$ws.Range("J91").Value = 1858.2858
$ws.Range("L91").Value = 1858.2858
$ws.Range("N91").Value = -4666.2858


$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 21450
$ws.Range("J82").Value = 38900
$ws.Range("L82").Value = 38900
$ws.Range("N82").Value = -39666

# Row 85
$ws.Range("H85").Value = 21450
$ws.Range("J85").Value = 38900
$ws.Range("L85").Value = 38900
$ws.Range("N85").Value = -41552

# Row 99
$ws.Range("H99").Value = 3722.9524
$ws.Range("I99").Value = 2541.8276
$ws.Range("K99").Value = 2541.8276
$ws.Range("M99").Value = -1043.8276


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4109.59
$ws.Range("I31").Value = 2990.1667
$ws.Range("J31").Value = 7841
$ws.Range("K31").Value = 2990.1667
$ws.Range("L31").Value = 7841
$ws.Range("M31").Value = -2695.1667
$ws.Range("N31").Value = -8431

# Row 34
$ws.Range("H34").Value = 4109.59
$ws.Range("I34").Value = 2990.1667
$ws.Range("J34").Value = 7841
$ws.Range("K34").Value = 2990.1667
$ws.Range("L34").Value = 7841
$ws.Range("M34").Value = -2788.1667
$ws.Range("N34").Value = -8245

# Row 41
$ws.Range("H41").Value = 13946.714
$ws.Range("I41").Value = 7931.75
$ws.Range("J41").Value = 21966.666
$ws.Range("K41").Value = 7931.75
$ws.Range("L41").Value = 21966.666
$ws.Range("M41").Value = -7503.75
$ws.Range("N41").Value = -22822.666

# Row 51
$ws.Range("H51").Value = 30116.334
$ws.Range("J51").Value = 30030
$ws.Range("L51").Value = 30030
$ws.Range("N51").Value = -31502

# Row 59
$ws.Range("H59").Value = 46396.668
$ws.Range("J59").Value = 46396.668
$ws.Range("L59").Value = 46396.668
$ws.Range("N59").Value = -48686.668

# Row 60
$ws.Range("H60").Value = 24387.5
$ws.Range("J60").Value = 29266.666
$ws.Range("L60").Value = 29266.666
$ws.Range("N60").Value = -30288.666

# Row 61
$ws.Range("H61").Value = 30116.334
$ws.Range("J61").Value = 30030
$ws.Range("L61").Value = 30030
$ws.Range("N61").Value = -30726

# Row 68
$ws.Range("H68").Value = 44496.668
$ws.Range("J68").Value = 44496.668
$ws.Range("L68").Value = 44496.668
$ws.Range("N68").Value = -45994.668

# Row 71
$ws.Range("H71").Value = 44496.668
$ws.Range("J71").Value = 44496.668
$ws.Range("L71").Value = 133490.004
$ws.Range("N71").Value = -140978.004

# Row 134
$ws.Range("H134").Value = 3850.3845
$ws.Range("I134").Value = 2038.375
$ws.Range("J134").Value = 6749.6
$ws.Range("K134").Value = 6115.125
$ws.Range("L134").Value = 20248.8
$ws.Range("M134").Value = -3580.125
$ws.Range("N134").Value = -25318.8


$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 6632047.5
$ws.Range("I11").Value = 7875369
$ws.Range("K11").Value = 23626107
$ws.Range("M11").Value = -23625967

# Row 14
$ws.Range("H14").Value = 224.77777
$ws.Range("I14").Value = 224.77777
$ws.Range("K14").Value = 674.33331
$ws.Range("M14").Value = -501.33331

# Row 38
$ws.Range("H38").Value = 639
$ws.Range("J38").Value = 1144
$ws.Range("L38").Value = 3432
$ws.Range("N38").Value = -4126

# Row 107
$ws.Range("H107").Value = 1283.75
$ws.Range("I107").Value = 1283.75
$ws.Range("K107").Value = 3851.25
$ws.Range("M107").Value = -1931.25

# Row 134
$ws.Range("H134").Value = 1310.125
$ws.Range("I134").Value = 1310.125
$ws.Range("K134").Value = 3930.375
$ws.Range("M134").Value = 1139.625

# Row 140
$ws.Range("H140").Value = 746527.0600000001
$ws.Range("I140").Value = 1432.6
$ws.Range("K140").Value = 4297.799999999999
$ws.Range("M140").Value = 882.2000000000007


$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4136.222
$ws.Range("I80").Value = 3681.125
$ws.Range("K80").Value = 3681.125
$ws.Range("M80").Value = -2683.125

# Row 83
$ws.Range("H83").Value = 4136.222
$ws.Range("I83").Value = 3681.125
$ws.Range("K83").Value = 18405.625
$ws.Range("M83").Value = -13413.625

# Row 93
$ws.Range("H93").Value = 36139
$ws.Range("J93").Value = 36139
$ws.Range("L93").Value = 36139
$ws.Range("N93").Value = -39883

# Row 97
$ws.Range("H97").Value = 7774.75
$ws.Range("I97").Value = 1639.8
$ws.Range("J97").Value = 17999.666
$ws.Range("K97").Value = 1639.8
$ws.Range("L97").Value = 17999.666
$ws.Range("M97").Value = -1143.8
$ws.Range("N97").Value = -18991.666

# Row 122
$ws.Range("H122").Value = 2169.8
$ws.Range("I122").Value = 2169.8
$ws.Range("K122").Value = 6509.400000000001
$ws.Range("M122").Value = -4059.400000000001

# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()


$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 8561.267
$ws.Range("I68").Value = 9060.25
$ws.Range("J68").Value = 7991
$ws.Range("K68").Value = 9060.25
$ws.Range("L68").Value = 7991
$ws.Range("M68").Value = -8311.25
$ws.Range("N68").Value = -9489

# Row 71
$ws.Range("H71").Value = 8561.267
$ws.Range("I71").Value = 9060.25
$ws.Range("J71").Value = 7991
$ws.Range("K71").Value = 45301.25
$ws.Range("L71").Value = 39955
$ws.Range("M71").Value = -41557.25
$ws.Range("N71").Value = -47443

# Row 93
$ws.Range("H93").Value = 15442.8
$ws.Range("J93").Value = 55200
$ws.Range("L93").Value = 55200
$ws.Range("N93").Value = -57696

# Row 132
$ws.Range("H132").Value = 10732.066
$ws.Range("I132").Value = 12670.2
$ws.Range("K132").Value = 38010.60000000001
$ws.Range("M132").Value = -35480.60000000001


$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 21160
$ws.Range("J54").Value = 29900
$ws.Range("L54").Value = 29900
$ws.Range("N54").Value = -30940

# Row 81
$ws.Range("H81").Value = 4225.2856
$ws.Range("J81").Value = 5638.25
$ws.Range("L81").Value = 11276.5
$ws.Range("N81").Value = -13398.5

# Row 84
$ws.Range("H84").Value = 4225.2856
$ws.Range("J84").Value = 5638.25
$ws.Range("L84").Value = 56382.5
$ws.Range("N84").Value = -66990.5

# Row 100
$ws.Range("H100").Value = 1156.8572
$ws.Range("I100").Value = 1216.3334
$ws.Range("J100").Value = 800
$ws.Range("K100").Value = 2432.6668
$ws.Range("L100").Value = 1600
$ws.Range("M100").Value = -1891.6668
$ws.Range("N100").Value = -2682

# Row 126
$ws.Range("H126").Value = 5356.115
$ws.Range("I126").Value = 5437.875
$ws.Range("K126").Value = 16313.625
$ws.Range("M126").Value = -13843.625
